$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Dentist
$ws.Range("A2").Value = "Dentist"
$ws.Range("B2").Value = 40000
$ws.Range("C2").Value = 46089.22928240741

# Row 3: Stipend
$ws.Range("A3").Value = "Stipend"
$ws.Range("B3").Value = 20000
$ws.Range("C3").Value = 46082.22928240741

# Row 4: Food Delivery (new row)
$ws.Range("A4").Value = "Food Delivery"
$ws.Range("B4").Value = 5000
$ws.Range("C4").Value = 46065.22928240741

# Row 5: Stocks (new row)
$ws.Range("A5").Value = "Stocks"
$ws.Range("B5").Value = 45000
$ws.Range("C5").Value = 46057.22928240741

# Row 6: Salary (new row)
$ws.Range("A6").Value = "Salary"
$ws.Range("B6").Value = 30000
$ws.Range("C6").Value = 46055.22928240741

# Row 7: Freelance (new row, moved from row 2)
$ws.Range("A7").Value = "Freelance"
$ws.Range("B7").Value = 10000
$ws.Range("C7").Value = 46054.22928240741

# Propagate the date number format (same as C3, numFmtId 14 => m/d/yyyy)
# from an already-formatted date cell onto the newly added date cells,
# so they share the same style index instead of minting a new one.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4:C7").PasteSpecial(-4122) | Out-Null
